$wb = $excel.ActiveWorkbook

# 1. Insert the new "Primary caregivers only" sheet after "parent 3 waves"
$afterSheet = $wb.Worksheets.Item("parent 3 waves")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "Primary caregivers only"

$ids = @("101_1", "102_1", "103_1", "104_1", "105_1", "106_1", "107_1", "108_1", "109_1", "110_1", "111_1", "112_1", "113_1", "114_1", "201_1", "202_1", "203_1", "204_1", "205_2", "206_1", "207_1", "208_1", "209_1", "210_1", "211_1", "212_1", "214_1", "215_1", "216_1", "302_1", "304_1", "305_1", "306_1", "307_1", "308_1", "309_1", "310_1", "311_1", "312_1", "313_1", "314_1", "315_2", "317_1", "318_1", "401_1", "403_1", "404_1", "405_2", "406_1", "407_1", "408_1", "410_1", "411_1", "412_1", "413_1", "414_1", "415_1", "416_1", "417_1", "418_1", "419_1", "420_1", "501_1", "502_1", "503_1", "504_1", "505_1", "506_1", "507_4", "508_1", "509_1", "510_1", "511_1", "513_1", "601_1", "602_1", "603_1", "604_1", "605_1", "606_1", "607_1", "608_1", "610_1", "611_1", "612_2", "613_1", "614_2", "615_1", "617_1", "618_1", "619_1", "620_1", "621_1", "622_1", "624_1")

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 1
    $newSheet.Cells.Item($row, 1).Value = $ids[$i]
}
$newSheet.Cells.Item(97, 1).Value = "n = 95"

Write-Output "done"
